$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text edits (before column insert, while columns are at their original positions) ---

# R32: remove the "Sponsored by..." prefix, keep the Yvonne Welbon intro note
$ws.Cells.Item(32, 18).Value = 'With introductory remarks by the director Yvonne Welbon.'

# H58: remove the trailing "Preceded by ..." sentence from the capsule description
$ws.Cells.Item(58, 8).Value = '_An Open Window_ is a film about mental illness among Dakar’s street wanderers. After shooting, Sylla fell ill, “seeing what others don’t see—the dislocated eye, the ancientness of the glass bubble, a sky descended too low, a horizon far too near. I was living the experience from the inside.” Fronza Wood’s _Killing Time_ is a dark comedy about finding the right outfit to commit suicide in. The New Yorker calls it “among the best short films ever made." After becoming the first Black woman in the New York camera operators'' union, Jessie Maple turned to directing. In _Twice As Nice_, twin basketball players compete to become the first woman in the “MBA.” The film stars real-life powerhouses Pamela and Paula McGhee, who led USC to back-to-back NCAA championships in the 1980s, and Cynthia Cooper-Dyke.'

# R58: remove the "Sponsored by..." phrase from the TWICE AS NICE note
$ws.Cells.Item(58, 18).Value = 'TWICE AS NICE initial screening (on 4/11/22) cancelled due to booking issues. Introduction by Danielle Scruggs'

# --- Clear the redundant "Sponsored by South Side Projection and Arts + Public Life" /
#     "Part of a screening of MONANGAMBEE..." public notes that are now consolidated
#     into the series-level note on the Series Info sheet ---
$ws.Cells.Item(4, 18).ClearContents()
$ws.Cells.Item(10, 18).ClearContents()
$ws.Cells.Item(24, 18).ClearContents()
$ws.Cells.Item(40, 18).ClearContents()
$ws.Cells.Item(49, 18).ClearContents()
$ws.Cells.Item(66, 18).ClearContents()

# --- Insert a new "ticketing url" column between "public notes" (R) and "slot" (old S, now T) ---
$newColWidth = $ws.Columns.Item(18).ColumnWidth
$ws.Columns.Item(19).Insert()
$ws.Columns.Item(19).ColumnWidth = $newColWidth
$ws.Cells.Item(1, 19).Value = "ticketing url"
